# "invalid_emails.xlsx" - create excel file without emailid customers
#
# The "Invalid Emails" sheet used to carry blank placeholders in column A
# (Name) for every customer row. Populate those rows with the customers'
# id numbers instead, leaving column B (Email) untouched/blank, and
# without disturbing the existing cell formatting.
#
# The id values look numeric, so a plain `.Value = "123..."` assignment
# would make Excel coerce them into real numbers (losing the shared
# string / text representation and any leading zeros). To keep them as
# literal text we stage each value in a scratch cell formatted as Text,
# copy it, and use PasteSpecial(xlPasteValues) so only the value moves
# into the target cell - the destination keeps its original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids = @(
    "7771080120987",
    "6750080045432",
    "1117080008888",
    "1117080088888",
    "1135088888888",
    "1153080077777",
    "2020080166666",
    "2020080135555",
    "2039080022222"
)

$xlPasteValues = -4163

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $scratch.Value = $ids[$i]
    $scratch.Copy()
    $ws.Range("A$row").PasteSpecial($xlPasteValues)
}

$scratch.Clear()
$excel.CutCopyMode = $false
